$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2412667512.985902
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 281020871479559.7
$ws.Range("G2").Value = 599.9793834530285
$ws.Range("H2").Value = 2412667512.985902

$ws.Range("C3").Value = 734245007.6737363
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 109666514042594.4
$ws.Range("G3").Value = 599.9839238902063
$ws.Range("H3").Value = 734245007.6737363

$ws.Range("C4").Value = 2476191769.954624
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 284066152935168.8
$ws.Range("G4").Value = 599.9790672944112
$ws.Range("H4").Value = 2476191769.954624

$ws.Range("C5").Value = 2267867024.890732
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 276589522591358.8
$ws.Range("G5").Value = 599.9803111714041
$ws.Range("H5").Value = 2267867024.890732

$ws.Range("C6").Value = 684729526.2772753
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 99363076310209.3
$ws.Range("G6").Value = 599.98345356569
$ws.Range("H6").Value = 684729526.2772753

$ws.Range("C7").Value = 88451370.87369135
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 13721556274646.81
$ws.Range("G7").Value = 599.984522831602
$ws.Range("H7").Value = 88451370.87369135

$ws.Range("C8").Value = 448083662.4119719
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 68321158008464.47
$ws.Range("G8").Value = 599.9842526323962
$ws.Range("H8").Value = 448083662.4119719

$ws.Range("C9").Value = 1410829212.825469
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 192719448351730.3
$ws.Range("G9").Value = 599.982422212462
$ws.Range("H9").Value = 1410829212.825469

$ws.Range("C10").Value = 11.26879859863915
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 7166599.296977292
$ws.Range("G10").Value = 599.9960259751688
$ws.Range("H10").Value = 11.26879859863915

$ws.Range("C11").Value = 80448460319.36288
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 6442078993411337
$ws.Range("G11").Value = 599.9700017775917
$ws.Range("H11").Value = 80448460319.36288

$ws.Range("C12").Value = 26997107502.8449
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 2651756066357992
$ws.Range("G12").Value = 599.9755477184607
$ws.Range("H12").Value = 26997107502.8449

$ws.Range("C13").Value = 82260972169.52925
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 6550165311036239
$ws.Range("G13").Value = 599.9698318209671
$ws.Range("H13").Value = 82260972169.52925

$ws.Range("C14").Value = 76853197667.61612
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 6443047166599908
$ws.Range("G14").Value = 599.9713481630326
$ws.Range("H14").Value = 76853197667.61612

$ws.Range("C15").Value = 24167913345.41976
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 2522363449560980
$ws.Range("G15").Value = 599.9769876669343
$ws.Range("H15").Value = 24167913345.41976

$ws.Range("C16").Value = 3220626050.680063
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 337088300014996.2
$ws.Range("G16").Value = 599.977054049558
$ws.Range("H16").Value = 3220626050.680063

$ws.Range("C17").Value = 16236655844.35864
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1692527030221570
$ws.Range("G17").Value = 599.9769593893565
$ws.Range("H17").Value = 16236655844.35864

$ws.Range("C18").Value = 50817480527.65192
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 4565004790096096
$ws.Range("G18").Value = 599.9732629300808
$ws.Range("H18").Value = 50817480527.65192

$ws.Range("C19").Value = 594.0351617115377
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 156848507.1189587
$ws.Range("G19").Value = 599.9908988153697
$ws.Range("H19").Value = 594.0351617115377
